# Updates cryptos list figures (price / 1h volume change, and a few
# reordered coin rows) to match the latest scrape.
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel stores them as text (matching the sheet's existing text format)
# instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.428.97"
$ws.Range("E2").Value = "  -1.63%  "

$ws.Range("D3").Value = "3.840.73"
$ws.Range("E3").Value = "  -1.07%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'601.37"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D6").Value = "'169.84"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "3.840.25"
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").Value = "'6.49"
$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("D14").Value = "'37.21"
$ws.Range("E14").Value = "  -2.85%  "

$ws.Range("D15").Value = "4.485.50"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").Value = "3.841.17"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "68.515.20"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").Value = "'18.50"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("D21").Value = "'11.10"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").Value = "'470.41"
$ws.Range("E22").Value = "  -4.04%  "

$ws.Range("D23").Value = "'0.740"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "'0.0000160"
$ws.Range("E24").Value = "  -4.20%  "

$ws.Range("D25").Value = "'83.32"
$ws.Range("E25").Value = "  -2.31%  "

$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -2.65%  "

$ws.Range("D27").Value = "'12.21"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  -1.34%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("D31").Value = "3.991.36"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").Value = "'7.69"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").Value = "'31.66"
$ws.Range("E33").Value = "  -0.84%  "

$ws.Range("E34").Value = "  -4.45%  "

$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").Value = "3.805.91"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("E37").Value = "  -1.80%  "

$ws.Range("E38").Value = "  +11.14%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'5.97"
$ws.Range("E39").Value = "  -2.51%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'1.02"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").Value = "'2.00"
$ws.Range("E44").Value = "  -5.72%  "

$ws.Range("D45").Value = "'8.79"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").Value = "'420.66"
$ws.Range("E46").Value = "  -4.02%  "

$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "'0.000292"
$ws.Range("E48").Value = "  +7.77%  "

$ws.Range("D49").Value = "'47.16"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0361"
$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'141.71"
$ws.Range("E51").Value = "  -1.31%  "

